# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed "K" (column G) values after switching the source metric from
# Strike# to K and recalculating std/mean based s_vals for each save row.
$newK = @{
    2  = 4
    3  = 0
    4  = 4
    5  = 1
    6  = 4
    7  = 2
    8  = 2
    9  = 2
    10 = 3
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 3
    17 = 0
    19 = 0
    20 = 2
    21 = 2
    22 = 2
    23 = 0
    24 = 1
    25 = 1
    26 = 7
    27 = 1
    28 = 3
    29 = 1
    30 = 1
    31 = 0
    32 = 2
    33 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
